# Ajout draft mapping f595a2bd5e53be80aa00972cfd76eee4a5f7087b
#
# - Bump the IG's "Date" metadata value.
# - Add a new mapping column ("Mapping: Spécification métier vers
#   l'extension ROR MetaComment") to the Elements data-dictionary sheet,
#   filled in only for the Extension.value[x] row.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date property -----------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- Elements sheet: add the new mapping column (AL) -----------------------
$wsElem = $wb.Worksheets.Item("Elements")

# Header cell (row 1), cloning the formatting of the preceding mapping
# column header so the new column matches the existing header style.
$wsElem.Range("AK1").Copy()
$wsElem.Range("AL1").PasteSpecial(-4122)
$wsElem.Range("AL1").Value = "Mapping: Spécification métier vers l'extension ROR MetaComment"

# Data cells (rows 2-6), cloning the formatting of the preceding mapping
# column's data style.
$wsElem.Range("AK2").Copy()
$wsElem.Range("AL2:AL6").PasteSpecial(-4122)

# Only the Extension.value[x] row (row 6) carries a mapping value; the rest
# of the column stays blank like the rest of the "Mapping: ..." columns.
$wsElem.Range("AL6").Value = "metadonnee.commentaire"

# Match the new column's width to the rest of the mapping columns.
$wsElem.Columns.Item(38).ColumnWidth = 69.1
